# Update simulation result values for pl_mw.xlsx ("case with 380 kV done")
# Rows 2-25 correspond to line indices 0-23; columns B:F,H,J,K hold the
# recomputed per-line results (G, I, L, M, N, O stay at 0, unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8699934474986151
$ws.Range("C2").Value = 0.09993956321301312
$ws.Range("D2").Value = 0.01606691811157646
$ws.Range("E2").Value = 0.0903083688686479
$ws.Range("F2").Value = 5.46564429882747
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.2594876274752238
$ws.Range("K2").Value = 0.8912979301010182
$ws.Range("B3").Value = 0.8491418512581674
$ws.Range("C3").Value = 0.09828364450071092
$ws.Range("D3").Value = 0.01618851453865489
$ws.Range("E3").Value = 0.09047103471660378
$ws.Range("F3").Value = 5.267556758085618
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2525470512218817
$ws.Range("K3").Value = 0.8712170596343753
$ws.Range("B4").Value = 0.837153390624195
$ws.Range("C4").Value = 0.09740540495499772
$ws.Range("D4").Value = 0.01631623572668417
$ws.Range("E4").Value = 0.09065840970483663
$ws.Range("F4").Value = 5.14657211899015
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.2484037763693294
$ws.Range("K4").Value = 0.8597948119226544
$ws.Range("B5").Value = 0.8324722751296179
$ws.Range("C5").Value = 0.09708217566191024
$ws.Range("D5").Value = 0.01638145929879187
$ws.Range("E5").Value = 0.09075674096193609
$ws.Range("F5").Value = 5.097426765661993
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.2467449465623304
$ws.Range("K5").Value = 0.8553676606780414
$ws.Range("B6").Value = 0.8317073021636077
$ws.Range("C6").Value = 0.09703059227217636
$ws.Range("D6").Value = 0.01639308057734468
$ws.Range("E6").Value = 0.0907743952135931
$ws.Range("F6").Value = 5.089275559796164
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2464712818516972
$ws.Range("K6").Value = 0.8546462567322806
$ws.Range("B7").Value = 0.8370894330472254
$ws.Range("C7").Value = 0.09740090565323101
$ws.Range("D7").Value = 0.01631706222480389
$ws.Range("E7").Value = 0.09065964689702
$ws.Range("F7").Value = 5.145908698422176
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.2483812852031662
$ws.Range("K7").Value = 0.8597341855507636
$ws.Range("B8").Value = 0.8626344308600551
$ws.Range("C8").Value = 0.09933975414928398
$ws.Range("D8").Value = 0.01609773145055371
$ws.Range("E8").Value = 0.09034627862548561
$ws.Range("F8").Value = 5.397207680379722
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.2570698762937838
$ws.Range("K8").Value = 0.8841852675648738
$ws.Range("B9").Value = 0.9192199988748371
$ws.Range("C9").Value = 0.1042489042386734
$ws.Range("D9").Value = 0.01609618782170941
$ws.Range("E9").Value = 0.09042750438568348
$ws.Range("F9").Value = 5.895328064028433
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2750545855681992
$ws.Range("K9").Value = 0.9393705433125774
$ws.Range("B10").Value = 0.964798845848037
$ws.Range("C10").Value = 0.1085426937193148
$ws.Range("D10").Value = 0.01636684744834582
$ws.Range("E10").Value = 0.09091373834473515
$ws.Range("F10").Value = 6.264898806538895
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2888577489137987
$ws.Range("K10").Value = 0.9843847241737933
$ws.Range("B11").Value = 0.9864146321399687
$ws.Range("C11").Value = 0.1106479643104024
$ws.Range("D11").Value = 0.01655128841516529
$ws.Range("E11").Value = 0.09122813009206254
$ws.Range("F11").Value = 6.433887887222625
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2952681591441291
$ws.Range("K11").Value = 1.005846633160559
$ws.Range("B12").Value = 0.994727555308657
$ws.Range("C12").Value = 0.1114672482124774
$ws.Range("D12").Value = 0.01663013704994398
$ws.Range("E12").Value = 0.09136062820130775
$ws.Range("F12").Value = 6.49801045641442
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2977147015141668
$ws.Range("K12").Value = 1.014116281974168
$ws.Range("B13").Value = 0.9929315368530069
$ws.Range("C13").Value = 0.1112898161570257
$ws.Range("D13").Value = 0.01661275208623891
$ws.Range("E13").Value = 0.09133149370796545
$ws.Range("F13").Value = 6.484194647867582
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2971869442582431
$ws.Range("K13").Value = 1.012328913499488
$ws.Range("B14").Value = 0.9870959820738392
$ws.Range("C14").Value = 0.1107149239084464
$ws.Range("D14").Value = 0.01655759378833466
$ws.Range("E14").Value = 0.09123876108847995
$ws.Range("F14").Value = 6.439160653637202
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.295469054633088
$ws.Range("K14").Value = 1.006524121277721
$ws.Range("B15").Value = 0.9835381609123033
$ws.Range("C15").Value = 0.1103656653198613
$ws.Range("D15").Value = 0.01652498603363028
$ws.Range("E15").Value = 0.0911837118782941
$ws.Range("F15").Value = 6.411593121777173
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2944192851358309
$ws.Range("K15").Value = 1.002987101509376
$ws.Range("B16").Value = 0.9634040056863284
$ws.Range("C16").Value = 0.1084081855217107
$ws.Range("D16").Value = 0.0163560436380763
$ws.Range("E16").Value = 0.09089507115394824
$ws.Range("F16").Value = 6.253872898396367
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2884414706377356
$ws.Range("K16").Value = 0.9830020284994987
$ws.Range("B17").Value = 0.9512786805746885
$ws.Range("C17").Value = 0.1072464118677487
$ws.Range("D17").Value = 0.01626824049620978
$ws.Range("E17").Value = 0.09074190028682594
$ws.Range("F17").Value = 6.157342912812396
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2848080270490954
$ws.Range("K17").Value = 0.9709946701609056
$ws.Range("B18").Value = 0.9443874646823929
$ws.Range("C18").Value = 0.1065924851836257
$ws.Range("D18").Value = 0.01622349853077054
$ws.Range("E18").Value = 0.09066257037889613
$ws.Range("F18").Value = 6.101902953871587
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2827305122229404
$ws.Range("K18").Value = 0.9641809884464863
$ws.Range("B19").Value = 0.9420684431736674
$ws.Range("C19").Value = 0.1063735257952771
$ws.Range("D19").Value = 0.01620933328221241
$ws.Range("E19").Value = 0.09063721560834992
$ws.Range("F19").Value = 6.083145787751249
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2820292163578131
$ws.Range("K19").Value = 0.9618898761129628
$ws.Range("B20").Value = 0.9525608520705759
$ws.Range("C20").Value = 0.107368603724467
$ws.Range("D20").Value = 0.01627698979154957
$ws.Range("E20").Value = 0.09075729763073426
$ws.Range("F20").Value = 6.16761021965155
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2851935342772833
$ws.Range("K20").Value = 0.9722632804257785
$ws.Range("B21").Value = 0.9888065601285234
$ws.Range("C21").Value = 0.1108831832359982
$ws.Range("D21").Value = 0.01657354922404153
$ws.Range("E21").Value = 0.09126563367233587
$ws.Range("F21").Value = 6.452384658763776
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2959731217367505
$ws.Range("K21").Value = 1.008225255150791
$ws.Range("B22").Value = 1.013238698459958
$ws.Range("C22").Value = 0.113308866412055
$ws.Range("D22").Value = 0.01681994346724736
$ws.Range("E22").Value = 0.09167624802801555
$ws.Range("F22").Value = 6.63926176364663
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.3031293833486188
$ws.Range("K22").Value = 1.032559451972588
$ws.Range("B23").Value = 1.000130541618233
$ws.Range("C23").Value = 0.1120023882722734
$ws.Range("D23").Value = 0.01668356509231472
$ws.Range("E23").Value = 0.09144990848888668
$ws.Range("F23").Value = 6.539450662569266
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2992997182290793
$ws.Range("K23").Value = 1.019495500402741
$ws.Range("B24").Value = 0.9519809338463006
$ws.Range("C24").Value = 0.1073133172176881
$ws.Range("D24").Value = 0.01627301637891776
$ws.Range("E24").Value = 0.0907503093025781
$ws.Range("F24").Value = 6.162968194887242
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2850192109636112
$ws.Range("K24").Value = 0.9716894629924866
$ws.Range("B25").Value = 0.9032122313034563
$ws.Range("C25").Value = 0.1028010900065368
$ws.Range("D25").Value = 0.01604964208171822
$ws.Range("E25").Value = 0.090330802844516
$ws.Range("F25").Value = 5.759966764273514
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2700866115771703
$ws.Range("K25").Value = 0.9236608224030647

Write-Host "Updated pl_mw.xlsx result values (rows 2-25, cols B-K) for the 380 kV case."
